# Trade #30 closed at 2026-02-17 04:15:27 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.62   # Current Capital
$summary.Range("B4").Value = 0.62      # Total P&L $
$summary.Range("B5").Value = 0.41      # Total P&L %
$summary.Range("B6").Value = 30        # Total Trades
$summary.Range("B7").Value = 14        # Winning Trades
$summary.Range("B9").Value = 46.67     # Win Rate %

# ---------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.62     # Capital
$status.Range("D4").Value = 30         # Trades
$status.Range("E4").Value = 0.62       # P&L $
$status.Range("F4").Value = 0.62       # P&L %
$status.Range("G4").Value = 46.67      # Win Rate %

# ---------------------------------------------------------------
# Append the new trade row (#30) to both "All Trades" and
# "MarketMaking" sheets - they mirror each other.
# ---------------------------------------------------------------
$newRow = 31

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item($newRow, 1).Value = 30

    # Column B holds a literal "yyyy-mm-dd" looking string, not an
    # actual date - force text storage so it isn't auto-converted to
    # a date serial by the smart-entry heuristic, then strip the
    # number-format override so the cell keeps the workbook's default
    # (unstyled) look, matching every other row in the sheet.
    $dateCell = $ws.Cells.Item($newRow, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026-02-17"
    $dateCell.ClearFormats()

    $ws.Cells.Item($newRow, 3).Value = "04:15:00"
    $ws.Cells.Item($newRow, 4).Value = "MarketMaking"
    $ws.Cells.Item($newRow, 5).Value = "UP"
    $ws.Cells.Item($newRow, 6).Value = 0.8
    $ws.Cells.Item($newRow, 7).Value = 0.99
    $ws.Cells.Item($newRow, 8).Value = "CLOSED"
    $ws.Cells.Item($newRow, 9).Value = 23.75
    $ws.Cells.Item($newRow, 10).Value = 0.19
    $ws.Cells.Item($newRow, 11).Value = 100.62
    $ws.Cells.Item($newRow, 12).Value = 0
    $ws.Cells.Item($newRow, 13).Value = 0
    $ws.Cells.Item($newRow, 14).Value = 0.6
    $ws.Cells.Item($newRow, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($newRow, 16).Value = "early_exit"
    $ws.Cells.Item($newRow, 17).Value = 0.44
}
